# Auto-generated edit script: Add data for 2023-09-06
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("E2").Value = 50
$ws.Range("C3").Value = 63
$ws.Range("G3").Value = 89
$ws.Range("I3").Value = 144
$ws.Range("B6").Value = 271
$ws.Range("C6").Value = 343
$ws.Range("D6").Value = 304
$ws.Range("E6").Value = 305
$ws.Range("F6").Value = 390
$ws.Range("H6").Value = 310
$ws.Range("I6").Value = 378
$ws.Range("J6").Value = 287
$ws.Range("B7").Value = 371
$ws.Range("C7").Value = 461
$ws.Range("D7").Value = 471
$ws.Range("E7").Value = 462
$ws.Range("F7").Value = 551
$ws.Range("G7").Value = 507
$ws.Range("H7").Value = 492
$ws.Range("I7").Value = 625
$ws.Range("J7").Value = 536

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("F6").Value = 44
$ws.Range("F7").Value = 51

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("E2").Value = 3
$ws.Range("H6").Value = 14
$ws.Range("E7").Value = 28
$ws.Range("H7").Value = 26

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J5").Value = 12
$ws.Range("J6").Value = 17

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 10

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("C3").Value = 1
$ws.Range("C6").Value = 9

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("C6").Value = 35
$ws.Range("C7").Value = 39

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J5").Value = 17
$ws.Range("I23").Value = 6
$ws.Range("G27").Value = 8
$ws.Range("C28").Value = 39
$ws.Range("F32").Value = 51
$ws.Range("E36").Value = 28
$ws.Range("H36").Value = 26
$ws.Range("J50").Value = 10
$ws.Range("E53").Value = 57
$ws.Range("I53").Value = 102
$ws.Range("J53").Value = 85
$ws.Range("J61").Value = 2
$ws.Range("G65").Value = 16
$ws.Range("B75").Value = 2
$ws.Range("E76").Value = 15
$ws.Range("F77").Value = 13
$ws.Range("C86").Value = 9
$ws.Range("D90").Value = 3
$ws.Range("B98").Value = 371
$ws.Range("C98").Value = 461
$ws.Range("D98").Value = 471
$ws.Range("E98").Value = 462
$ws.Range("F98").Value = 551
$ws.Range("G98").Value = 507
$ws.Range("H98").Value = 492
$ws.Range("I98").Value = 625
$ws.Range("J98").Value = 536

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 25
$ws.Range("E6").Value = 46
$ws.Range("I6").Value = 64
$ws.Range("J6").Value = 43
$ws.Range("E7").Value = 57
$ws.Range("I7").Value = 102
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("E5").Value = 7
$ws.Range("E6").Value = 15

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("G3").Value = 3
$ws.Range("G6").Value = 16

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I5").Value = 3
$ws.Range("I6").Value = 6

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 13

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("G3").Value = 4
$ws.Range("G5").Value = 8

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 3
